$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Drop the "_GoBack" bookmark that currently sits in the title
#    paragraph ("A3 Cross-Side-Scripting"). It gets re-created further
#    down, at the end of the bullet list, once the new content has been
#    written out.
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------
# 2) Locate the bullet paragraph that currently reads:
#    "Durch täuschend echt wirkende Benutzeranmeldeseiten, welche aber
#    bösartiger Herkunft sind, obwohl sie eine autoritäre URL aufweisen"
#    and replace it (and its trailing paragraph mark) with the new set
#    of bullets describing the XSS-based phishing attack in more detail.
# ---------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Durch*autorit*re URL aufweisen*") {
        $target = $para
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the 'Durch taeuschend echt wirkende ...' bullet paragraph"
}

$rng = $d.Range($target.Range.Start, $target.Range.End)

$newXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr>
      <w:ilvl w:val="2"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Finden einer XSS L&#252;cke</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> durch Eingabe von HTML Anweisungen</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Wenn die Anweisung ausgef&#252;hrt wird </w:t>
  </w:r>
  <w:r>
    <w:sym w:font="Wingdings" w:char="F0E0"/>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> L&#252;cke gegeben</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr>
      <w:ilvl w:val="2"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Seite mit eigenen Skripten innerhalb der Website erstellen</w:t>
  </w:r>
  <w:r>
    <w:t>, welche einer Anmeldeseite nachempfunden ist.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr>
      <w:ilvl w:val="2"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>An potentielle Opfer die URL zur b&#246;sartigen Seite schicken</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Listenabsatz"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Unter dem Vorwand, dass das Opfer sich einloggen m&#252;sse um ein Problem mit seinem Konto zu l&#246;sen.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($newXml)
